$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '64.536.82'
$ws.Range('E2').Value = '  -0.98%  '
Set-TextValue $ws.Range('D3') '3.439.91'
$ws.Range('E3').Value = '  -1.24%  '
$ws.Range('E4').Value = '  -0.04%  '
Set-TextValue $ws.Range('D5') '572.91'
$ws.Range('E5').Value = '  -1.17%  '
$ws.Range('E6').Value = '  -2.27%  '
$ws.Range('E7').Value = '  -0.01%  '
Set-TextValue $ws.Range('D8') '3.439.32'
$ws.Range('E8').Value = '  -1.28%  '
Set-TextValue $ws.Range('D9') '0.574'
$ws.Range('E9').Value = '  -6.31%  '
$ws.Range('E10').Value = '  -0.89%  '
$ws.Range('E11').Value = '  -3.18%  '
Set-TextValue $ws.Range('D12') '0.438'
$ws.Range('E12').Value = '  -2.28%  '
Set-TextValue $ws.Range('D13') '4.033.83'
$ws.Range('E13').Value = '  -1.41%  '
$ws.Range('E14').Value = '  -0.35%  '
$ws.Range('E15').Value = '  -4.00%  '
Set-TextValue $ws.Range('D16') '0.0000175'
$ws.Range('E16').Value = '  -9.52%  '
Set-TextValue $ws.Range('D17') '64.635.10'
$ws.Range('E17').Value = '  -0.91%  '
Set-TextValue $ws.Range('D18') '3.434.41'
$ws.Range('E18').Value = '  -1.43%  '
$ws.Range('E19').Value = '  -5.07%  '
Set-TextValue $ws.Range('D20') '13.71'
$ws.Range('E20').Value = '  -4.94%  '
Set-TextValue $ws.Range('D21') '377.93'
$ws.Range('E21').Value = '  -1.48%  '
$ws.Range('E22').Value = '  -3.68%  '
Set-TextValue $ws.Range('D24') '72.07'
$ws.Range('E24').Value = '  -0.90%  '
Set-TextValue $ws.Range('D25') '0.527'
$ws.Range('E25').Value = '  -4.70%  '
$ws.Range('E26').Value = '  -1.05%  '
Set-TextValue $ws.Range('D27') '9.90'
$ws.Range('E27').Value = '  -1.77%  '
$ws.Range('E28').Value = '  -0.38%  '
$ws.Range('E29').Value = '  +0.03%  '
Set-TextValue $ws.Range('D30') '1.44'
$ws.Range('E30').Value = '  -6.39%  '
$ws.Range('E31').Value = '  -2.88%  '
$ws.Range('E32').Value = '  -2.81%  '
Set-TextValue $ws.Range('D33') '23.20'
$ws.Range('E33').Value = '  -2.16%  '
Set-TextValue $ws.Range('D34') '6.98'
$ws.Range('E34').Value = '  -3.09%  '
$ws.Range('E35').Value = '  -4.48%  '
Set-TextValue $ws.Range('D36') '160.73'
$ws.Range('E36').Value = '  -0.80%  '
$ws.Range('E37').Value = '  -3.63%  '
Set-TextValue $ws.Range('D38') '2.882.80'
$ws.Range('E38').Value = '  -4.09%  '
Set-TextValue $ws.Range('D39') '0.0745'
$ws.Range('E39').Value = '  -4.75%  '
Set-TextValue $ws.Range('D40') '26.16'
$ws.Range('E40').Value = '  -3.02%  '
Set-TextValue $ws.Range('D41') '0.790'
$ws.Range('E41').Value = '  +0.99%  '
$ws.Range('E42').Value = '  -0.10%  '
$ws.Range('E43').Value = '  -1.81%  '
$ws.Range('E44').Value = '  -4.42%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D45') '0.0309'
$ws.Range('E45').Value = '  -4.26%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range('D46') '25.70'
$ws.Range('E46').Value = '  -1.11%  '
$ws.Range('E47').Value = '  +9.06%  '
Set-TextValue $ws.Range('D48') '320.35'
$ws.Range('E48').Value = '  +0.57%  '
$ws.Range('E49').Value = '  -3.97%  '
$ws.Range('E50').Value = '  -3.58%  '
Set-TextValue $ws.Range('D51') '0.841'
$ws.Range('E51').Value = '  -4.01%  '
